$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.075.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.77%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.268.61'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.01%  '

$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.50%  '

$ws.Range("E7").Value = '  -0.25%  '

$ws.Range("E8").Value = '  +2.56%  '

$ws.Range("E9").Value = '  -1.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.94%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.408'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.839.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.07%  '

$ws.Range("E13").Value = '  +1.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.47'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '68.096.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.47%  '

$ws.Range("E16").Value = '  -0.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.254.34'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.31'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '416.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.58%  '

$ws.Range("E21").Value = '  -0.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.38%  '

$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("E24").Value = '  -0.38%  '

$ws.Range("E25").Value = '  +0.72%  '

$ws.Range("E26").Value = '  +1.21%  '

$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.994'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.76%  '

$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.87'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.98%  '

$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '163.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.22%  '

$ws.Range("E36").Value = '  -0.98%  '

$ws.Range("E37").Value = '  +1.60%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.52%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.797'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.25%  '

$ws.Range("E40").Value = '  -1.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.34'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.48%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.644.11'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0676'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '337.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.77%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0273'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.83%  '

$ws.Range("E49").Value = '  +0.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.977'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.10%  '

$ws.Range("E51").Value = '  -0.51%  '
